$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 1398
$ws.Range("I6").Value = 1219.8889
$ws.Range("J6").Value = 2199.5
$ws.Range("K6").Value = 3659.6667
$ws.Range("L6").Value = 6598.5
$ws.Range("M6").Value = -3547.6667
$ws.Range("N6").Value = -6822.5
# Row 8
$ws.Range("H8").Value = 81.2
$ws.Range("I8").Value = 81.2
$ws.Range("K8").Value = 243.6
$ws.Range("M8").Value = -104.6
# Row 17
$ws.Range("H17").Value = 7953.8125
$ws.Range("J17").Value = 7953.8125
$ws.Range("L17").Value = 23861.4375
$ws.Range("N17").Value = -24197.4375
# Row 33
$ws.Range("H33").Value = 274.44644
$ws.Range("I33").Value = 197.875
$ws.Range("K33").Value = 197.875
$ws.Range("M33").Value = 31.125
# Row 38
$ws.Range("H38").Value = 717.8889
$ws.Range("I38").Value = 85.166664
$ws.Range("J38").Value = 1983.3334
$ws.Range("K38").Value = 255.499992
$ws.Range("L38").Value = 5950.0002
$ws.Range("M38").Value = 116.500008
$ws.Range("N38").Value = -6694.0002
# Row 70
$ws.Range("H70").Value = 1183.3334
$ws.Range("I70").Value = 1000
$ws.Range("J70").Value = 1244.4445
$ws.Range("K70").Value = 3000
$ws.Range("L70").Value = 3733.3335
$ws.Range("M70").Value = -2730
$ws.Range("N70").Value = -4273.333500000001
# Row 73
$ws.Range("H73").Value = 1183.3334
$ws.Range("I73").Value = 1000
$ws.Range("J73").Value = 1244.4445
$ws.Range("K73").Value = 3000
$ws.Range("L73").Value = 3733.3335
$ws.Range("M73").Value = -2064
$ws.Range("N73").Value = -5605.333500000001
# Row 112
$ws.Range("H112").Value = 1309
$ws.Range("J112").Value = 1337.1428
$ws.Range("L112").Value = 4011.4284
$ws.Range("N112").Value = -6227.428400000001
# Row 126
$ws.Range("H126").Value = 44963
$ws.Range("J126").Value = 44963
$ws.Range("L126").Value = 44963
$ws.Range("N126").Value = -54843
# Row 129
$ws.Range("H129").Value = 1264.8723
$ws.Range("I129").Value = 1691.625
$ws.Range("J129").Value = 1177.3334
$ws.Range("K129").Value = 5074.875
$ws.Range("L129").Value = 3532.0002
$ws.Range("M129").Value = -74.875
$ws.Range("N129").Value = -13532.0002
# Row 138
$ws.Range("H138").Value = 1727.8235
$ws.Range("I138").Value = 1475.091
$ws.Range("J138").Value = 1966.1143
$ws.Range("K138").Value = 4425.272999999999
$ws.Range("L138").Value = 5898.3429
$ws.Range("M138").Value = 714.7270000000008
$ws.Range("N138").Value = -16178.3429

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1470.5883
$ws.Range("I45").Value = 1430
$ws.Range("J45").Value = 1528.5714
$ws.Range("K45").Value = 1430
$ws.Range("L45").Value = 1528.5714
$ws.Range("M45").Value = -1053
$ws.Range("N45").Value = -2282.5714
# Row 46
$ws.Range("H46").Value = 8384
$ws.Range("J46").Value = 8384
$ws.Range("L46").Value = 8384
$ws.Range("N46").Value = -9022
# Row 102
$ws.Range("H102").Value = 111135900
$ws.Range("I102").Value = 166669000
$ws.Range("K102").Value = 166669000
$ws.Range("M102").Value = -166667378
# Row 122
$ws.Range("H122").Value = 2084.8918
$ws.Range("I122").Value = 2182.3462
$ws.Range("J122").Value = 1854.5454
$ws.Range("K122").Value = 6547.0386
$ws.Range("L122").Value = 5563.6362
$ws.Range("M122").Value = -4097.0386
$ws.Range("N122").Value = -10463.6362

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2209.1904
$ws.Range("I105").Value = 1881.0625
$ws.Range("J105").Value = 3259.2
$ws.Range("K105").Value = 1881.0625
$ws.Range("L105").Value = 3259.2
$ws.Range("M105").Value = -134.0625
$ws.Range("N105").Value = -6753.2
# Row 119
$ws.Range("H119").Value = 47753
$ws.Range("J119").Value = 47753
$ws.Range("L119").Value = 47753
$ws.Range("N119").Value = -57429

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1767.6129
$ws.Range("I58").Value = 1469.8148
$ws.Range("J58").Value = 3777.75
$ws.Range("K58").Value = 1469.8148
$ws.Range("L58").Value = 3777.75
$ws.Range("M58").Value = -1266.8148
$ws.Range("N58").Value = -4183.75
# Row 62
$ws.Range("H62").Value = 2793.889
$ws.Range("I62").Value = 2793.889
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2793.889
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2169.889
$ws.Range("N62").Value = $null
# Row 65
$ws.Range("H65").Value = 2793.889
$ws.Range("I65").Value = 2793.889
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 13969.445
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -10849.445
$ws.Range("N65").Value = $null
# Row 132
$ws.Range("H132").Value = 65845.55
$ws.Range("I132").Value = 1706.2941
$ws.Range("J132").Value = 283919
$ws.Range("K132").Value = 5118.8823
$ws.Range("L132").Value = 851757
$ws.Range("M132").Value = -2588.8823
$ws.Range("N132").Value = -856817
# Row 134
$ws.Range("H134").Value = 610458.7
$ws.Range("I134").Value = 1095.5883
$ws.Range("J134").Value = 2336987.5
$ws.Range("K134").Value = 3286.7649
$ws.Range("L134").Value = 7010962.5
$ws.Range("M134").Value = -751.7648999999997
$ws.Range("N134").Value = -7016032.5
# Row 136
$ws.Range("H136").Value = 1767.6129
$ws.Range("I136").Value = 1469.8148
$ws.Range("J136").Value = 3777.75
$ws.Range("K136").Value = 4409.4444
$ws.Range("L136").Value = 11333.25
$ws.Range("M136").Value = -1859.4444
$ws.Range("N136").Value = -16433.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 35714604
$ws.Range("I2").Value = 24.4
$ws.Range("J2").Value = 68182400
$ws.Range("K2").Value = 146.4
$ws.Range("L2").Value = 409094400
$ws.Range("M2").Value = -33.39999999999998
$ws.Range("N2").Value = -409094626
# Row 38
$ws.Range("H38").Value = 17928786
$ws.Range("I38").Value = 135.26315
$ws.Range("J38").Value = 55778160
$ws.Range("K38").Value = 405.78945
$ws.Range("L38").Value = 167334480
$ws.Range("M38").Value = -58.78944999999999
$ws.Range("N38").Value = -167335174

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2730.4
$ws.Range("I102").Value = 2036.25
$ws.Range("J102").Value = 5507
$ws.Range("K102").Value = 2036.25
$ws.Range("L102").Value = 5507
$ws.Range("M102").Value = -414.25
$ws.Range("N102").Value = -8751
# Row 126
$ws.Range("H126").Value = 31254872
$ws.Range("I126").Value = 50006710
$ws.Range("J126").Value = 1807.3334
$ws.Range("K126").Value = 150020130
$ws.Range("L126").Value = 5422.0002
$ws.Range("M126").Value = -150017660
$ws.Range("N126").Value = -10362.0002
# Row 130
$ws.Range("H130").Value = 46403.816
$ws.Range("J130").Value = 46403.816
$ws.Range("L130").Value = 46403.816
$ws.Range("N130").Value = -56443.816

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 1280.4
$ws.Range("I82").Value = 1400.1428
$ws.Range("J82").Value = 1001
$ws.Range("K82").Value = 1400.1428
$ws.Range("L82").Value = 1001
$ws.Range("M82").Value = -1039.1428
$ws.Range("N82").Value = -1723
# Row 85
$ws.Range("H85").Value = 1280.4
$ws.Range("I85").Value = 1400.1428
$ws.Range("J85").Value = 1001
$ws.Range("K85").Value = 1400.1428
$ws.Range("L85").Value = 1001
$ws.Range("M85").Value = -152.1428000000001
$ws.Range("N85").Value = -3497
# Row 100
$ws.Range("H100").Value = 3297.6667
$ws.Range("I100").Value = 2893
$ws.Range("J100").Value = 3500
$ws.Range("K100").Value = 2893
$ws.Range("L100").Value = 3500
$ws.Range("M100").Value = -2352
$ws.Range("N100").Value = -4582

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1000
$ws.Range("I96").Value = 1000
$ws.Range("J96").Value = 1000
$ws.Range("K96").Value = 1000
$ws.Range("L96").Value = 1000
$ws.Range("M96").Value = 373
$ws.Range("N96").Value = -3746
